$wb = $excel.ActiveWorkbook

# This workbook lists 江西 (Jiangxi) comic/anime conventions. Bilibili view
# counts (column F, "想去人数") and occasionally the minimum ticket price
# (column G, "最低票价") were refreshed, and one brand-new event
# (上饶·星河城市动漫文化节, 2024-11-16) was inserted into the date-sorted
# listing, which pushes the three rows that used to sit at 30-32 down to
# 31-33 (their content changes accordingly) and grows the sheet by one row
# (A1:I32 -> A1:I33).
#
# Both the "展览" (sheet 1) and "全部类型" (sheet 4) worksheets contain an
# identical copy of this table, so every edit is applied to both.

foreach ($sheetIndex in 1,4) {
  $ws = $wb.Worksheets.Item($sheetIndex)

  # --- Refreshed "want to go" counts (column F) ---
  $ws.Cells.Item(2,6).Value = 16
  $ws.Cells.Item(3,6).Value = 85
  $ws.Cells.Item(4,6).Value = 257
  $ws.Cells.Item(5,6).Value = 44
  $ws.Cells.Item(6,6).Value = 535
  $ws.Cells.Item(7,6).Value = 49
  $ws.Cells.Item(8,6).Value = 1981
  $ws.Cells.Item(9,6).Value = 65
  $ws.Cells.Item(10,6).Value = 90
  $ws.Cells.Item(11,6).Value = 4213
  $ws.Cells.Item(13,6).Value = 274
  $ws.Cells.Item(14,6).Value = 96
  $ws.Cells.Item(15,6).Value = 92
  $ws.Cells.Item(16,6).Value = 17
  $ws.Cells.Item(17,6).Value = 56
  $ws.Cells.Item(18,6).Value = 2941
  $ws.Cells.Item(19,6).Value = 53
  $ws.Cells.Item(20,6).Value = 413
  $ws.Cells.Item(23,6).Value = 63
  $ws.Cells.Item(24,6).Value = 14
  $ws.Cells.Item(25,6).Value = 66
  $ws.Cells.Item(27,6).Value = 6
  $ws.Cells.Item(28,6).Value = 43
  $ws.Cells.Item(29,6).Value = 192

  # Row 11 minimum ticket price (column G) was refreshed too
  $ws.Cells.Item(11,7).Value = 65

  # --- New event inserted -> rows 30-32 shift down to 31-33 ---
  # Row 33 (new): old row 32 content (南昌·云芽动漫音乐嘉年华·封茗囧菌内场票),
  # with its "want to go" count refreshed from 234 to 235.
  $ws.Cells.Item(33,1).Value = 32
  $ws.Cells.Item(33,2).Formula = "'2024-12-08"
  $ws.Cells.Item(33,3).Value = "南昌·云芽动漫音乐嘉年华·封茗囧菌内场票"
  $ws.Cells.Item(33,4).Value = "怀玉山大道1315号 南昌绿地国际博览中心"
  $ws.Cells.Item(33,5).Value = "2024.12.08 09:30-12.08 17:30"
  $ws.Cells.Item(33,6).Value = 235
  $ws.Cells.Item(33,7).Value = 128
  $ws.Cells.Item(33,8).Value = "https://show.bilibili.com/platform/detail.html?id=92134"
  $ws.Cells.Item(33,9).Value = "//i0.hdslb.com/bfs/openplatform/202409/eeFHJb3W1725328994111.jpeg"

  # Row 32 (new): old row 31 content (南昌·云芽动漫音乐嘉年华), with its
  # "want to go" count refreshed from 1638 to 1654.
  $ws.Cells.Item(32,1).Value = 31
  $ws.Cells.Item(32,2).Formula = "'2024-12-07"
  $ws.Cells.Item(32,3).Value = "南昌·云芽动漫音乐嘉年华"
  $ws.Cells.Item(32,4).Value = "怀玉山大道1315号 南昌绿地国际博览中心"
  $ws.Cells.Item(32,5).Value = "2024.12.07 09:00-12.08 18:00"
  $ws.Cells.Item(32,6).Value = 1654
  $ws.Cells.Item(32,7).Value = 69
  $ws.Cells.Item(32,8).Value = "https://show.bilibili.com/platform/detail.html?id=92144"
  $ws.Cells.Item(32,9).Value = "//i0.hdslb.com/bfs/openplatform/202409/2DwZA4qv1725706772865.png"

  # Row 31 (new): old row 30 content (南昌·CM04动漫游戏博览会), with its
  # "want to go" count refreshed from 319 to 358.
  $ws.Cells.Item(31,1).Value = 30
  $ws.Cells.Item(31,2).Formula = "'2024-11-16"
  $ws.Cells.Item(31,3).Value = "南昌·CM04动漫游戏博览会"
  $ws.Cells.Item(31,4).Value = "怀玉山大道1315号 南昌绿地国际博览中心"
  $ws.Cells.Item(31,5).Value = "2024.11.16 09:00-11.17 17:00"
  $ws.Cells.Item(31,6).Value = 358
  $ws.Cells.Item(31,7).Value = 65
  $ws.Cells.Item(31,8).Value = "https://show.bilibili.com/platform/detail.html?id=92378"
  $ws.Cells.Item(31,9).Value = "//i2.hdslb.com/bfs/openplatform/202409/N57Jfogr1725381095803.jpeg"

  # Row 30 (new event): 上饶·星河城市动漫文化节
  $ws.Cells.Item(30,1).Value = 29
  $ws.Cells.Item(30,2).Formula = "'2024-11-16"
  $ws.Cells.Item(30,3).Value = "上饶·星河城市动漫文化节"
  $ws.Cells.Item(30,4).Value = "春江北大道时光PARK内 博悦宴会艺术中心"
  $ws.Cells.Item(30,5).Value = "2024.11.16 10:00-11.16 17:00"
  $ws.Cells.Item(30,6).Value = 4
  $ws.Cells.Item(30,7).Value = 55
  $ws.Cells.Item(30,8).Value = "https://show.bilibili.com/platform/detail.html?id=92572"
  $ws.Cells.Item(30,9).Value = "//i1.hdslb.com/bfs/openplatform/202409/sb00Vz8W1726637703913.jpeg"
}
